$d = $word.ActiveDocument

# Step 1: Remove the paragraph referencing the Distributed Transaction
# Coordinator entirely (its text and its own paragraph mark).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Distributed Transaction Coordinator*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# Step 2: Re-seat the "_GoBack" bookmark (Word always keeps exactly one,
# tracking the location of the most recent edit) at the start of the
# paragraph that now follows directly -- "Locate openXDA in the list...".
$dest = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Locate openXDA in the list*") {
        $dest = $p
        break
    }
}
if ($dest -ne $null) {
    $r = $d.Range($dest.Range.Start, $dest.Range.Start)
    $d.Bookmarks.Add("_GoBack", $r)
}

Write-Output "done"
